$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.48215109446573
$ws.Range("C2").Value = 9.622915436724618
$ws.Range("D2").Value = 6.929506622363342
$ws.Range("F2").Value = 37.95862912282813
$ws.Range("G2").Value = 3.680162490946819
$ws.Range("I2").Value = 29.98556928058724
$ws.Range("L2").Value = 10.77904386553931
$ws.Range("B3").Value = 19.02328921379115
$ws.Range("C3").Value = 8.99480798270759
$ws.Range("D3").Value = 6.947657735997654
$ws.Range("F3").Value = 37.55550981691331
$ws.Range("G3").Value = 3.684471620375105
$ws.Range("I3").Value = 29.87936398640711
$ws.Range("L3").Value = 10.75831953067555
$ws.Range("B4").Value = 18.74230747152236
$ws.Range("C4").Value = 8.58541395486165
$ws.Range("D4").Value = 6.959745809985331
$ws.Range("F4").Value = 37.31756767505772
$ws.Range("G4").Value = 3.687251992537754
$ws.Range("I4").Value = 29.82146239301849
$ws.Range("L4").Value = 10.74797073680646
$ws.Range("B5").Value = 18.62818656717883
$ws.Range("C5").Value = 8.412599992245919
$ws.Range("D5").Value = 6.964907749490089
$ws.Range("F5").Value = 37.22310176012114
$ws.Range("G5").Value = 3.688418989765443
$ws.Range("I5").Value = 29.79970957796863
$ws.Range("L5").Value = 10.74435290624751
$ws.Range("B6").Value = 18.60926556011419
$ws.Range("C6").Value = 8.383542631097853
$ws.Range("D6").Value = 6.965779100733915
$ws.Range("F6").Value = 37.20756910750517
$ws.Range("G6").Value = 3.688614824742956
$ws.Range("I6").Value = 29.79620891458498
$ws.Range("L6").Value = 10.74378841566522
$ws.Range("B7").Value = 18.74076659822068
$ws.Range("C7").Value = 8.583107566882823
$ws.Range("D7").Value = 6.959814471996858
$ws.Range("F7").Value = 37.31628344891851
$ws.Range("G7").Value = 3.687267593332825
$ws.Range("I7").Value = 29.82116156236169
$ws.Range("L7").Value = 10.74791951643785
$ws.Range("B8").Value = 19.3238898911937
$ws.Range("C8").Value = 9.411260685762917
$ws.Range("D8").Value = 6.935568620895028
$ws.Range("F8").Value = 37.81770206449075
$ws.Range("G8").Value = 3.681620437259089
$ws.Range("I8").Value = 29.94743624901157
$ws.Range("L8").Value = 10.77140626675487
$ws.Range("B9").Value = 20.46538186172601
$ws.Range("C9").Value = 10.84767403863323
$ws.Range("D9").Value = 6.895561838619815
$ws.Range("F9").Value = 38.87272949663448
$ws.Range("G9").Value = 3.67160760041422
$ws.Range("I9").Value = 30.25274375864871
$ws.Range("L9").Value = 10.83621128714113
$ws.Range("B10").Value = 21.29256561129516
$ws.Range("C10").Value = 11.78971273841944
$ws.Range("D10").Value = 6.870841249272663
$ws.Range("F10").Value = 39.68561082477194
$ws.Range("G10").Value = 3.664889192929523
$ws.Range("I10").Value = 30.51165645308195
$ws.Range("L10").Value = 10.89508137823243
$ws.Range("B11").Value = 21.66439817853377
$ws.Range("C11").Value = 12.19391027428599
$ws.Range("D11").Value = 6.860626431017707
$ws.Range("F11").Value = 40.06226389608983
$ws.Range("G11").Value = 3.661969443879082
$ws.Range("I11").Value = 30.63678867230391
$ws.Range("L11").Value = 10.92426058029865
$ws.Range("B12").Value = 21.80440941610894
$ws.Range("C12").Value = 12.34348843794052
$ws.Range("D12").Value = 6.856907967675014
$ws.Range("F12").Value = 40.20576450777135
$ws.Range("G12").Value = 3.660883289874376
$ws.Range("I12").Value = 30.68521230356493
$ws.Range("L12").Value = 10.93565014805816
$ws.Range("B13").Value = 21.77429310397773
$ws.Range("C13").Value = 12.31142871581876
$ws.Range("D13").Value = 6.85770212551373
$ws.Range("F13").Value = 40.17482225709779
$ws.Range("G13").Value = 3.661116347830148
$ws.Range("I13").Value = 30.67473752285527
$ws.Range("L13").Value = 10.93318216101699
$ws.Range("B14").Value = 21.67593372927799
$ws.Range("C14").Value = 12.20628585671864
$ws.Range("D14").Value = 6.860317502881403
$ws.Range("F14").Value = 40.07405301220188
$ws.Range("G14").Value = 3.661879695468966
$ws.Range("I14").Value = 30.64075180760258
$ws.Range("L14").Value = 10.9251908219167
$ws.Range("B15").Value = 21.61557806765132
$ws.Range("C15").Value = 12.1414297961762
$ws.Range("D15").Value = 6.861939031435462
$ws.Range("F15").Value = 40.01243877293315
$ws.Range("G15").Value = 3.662349802309588
$ws.Range("I15").Value = 30.62006927414426
$ws.Range("L15").Value = 10.9203400246753
$ws.Range("B16").Value = 21.26816247712696
$ws.Range("C16").Value = 11.76280917167382
$ws.Range("D16").Value = 6.871529675871608
$ws.Range("F16").Value = 39.66112386295315
$ws.Range("G16").Value = 3.665082740321382
$ws.Range("I16").Value = 30.5036251780822
$ws.Range("L16").Value = 10.89322231096361
$ws.Range("B17").Value = 21.05378045652679
$ws.Range("C17").Value = 11.52431581879995
$ws.Range("D17").Value = 6.877678241897961
$ws.Range("F17").Value = 39.44727964064552
$ws.Range("G17").Value = 3.666794170872135
$ws.Range("I17").Value = 30.43406100732327
$ws.Range("L17").Value = 10.87719753876288
$ws.Range("B18").Value = 20.93006247230908
$ws.Range("C18").Value = 11.38484863233166
$ws.Range("D18").Value = 6.881311661174914
$ws.Range("F18").Value = 39.32493579365078
$ws.Range("G18").Value = 3.667791395449371
$ws.Range("I18").Value = 30.39474273508459
$ws.Range("L18").Value = 10.86820663527759
$ws.Range("B19").Value = 20.88810791251687
$ws.Range("C19").Value = 11.33723323815209
$ws.Range("D19").Value = 6.882558477732056
$ws.Range("F19").Value = 39.28362810642943
$ws.Range("G19").Value = 3.668131250536142
$ws.Range("I19").Value = 30.38154982608767
$ws.Range("L19").Value = 10.86520144856138
$ws.Range("B20").Value = 21.07664547884642
$ws.Range("C20").Value = 11.54994094301973
$ws.Range("D20").Value = 6.877013674529431
$ws.Range("F20").Value = 39.46997691689161
$ws.Range("G20").Value = 3.66661065670617
$ws.Range("I20").Value = 30.44139461498745
$ws.Range("L20").Value = 10.87888003453229
$ws.Range("B21").Value = 21.70484696437908
$ws.Range("C21").Value = 12.23726324126358
$ws.Range("D21").Value = 6.859545229583184
$ws.Range("F21").Value = 40.10362870367138
$ws.Range("G21").Value = 3.661654953843425
$ws.Range("I21").Value = 30.65070620062168
$ws.Range("L21").Value = 10.92752888498667
$ws.Range("B22").Value = 22.11072375551511
$ws.Range("C22").Value = 12.66617670135801
$ws.Range("D22").Value = 6.849001647980868
$ws.Range("F22").Value = 40.52277450182375
$ws.Range("G22").Value = 3.658529664854129
$ws.Range("I22").Value = 30.79354823688202
$ws.Range("L22").Value = 10.96130310503059
$ws.Range("B23").Value = 21.89457631602922
$ws.Range("C23").Value = 12.43910834211046
$ws.Range("D23").Value = 6.854548587004718
$ws.Range("F23").Value = 40.29864835985355
$ws.Range("G23").Value = 3.660187345263847
$ws.Range("I23").Value = 30.71676425280274
$ws.Range("L23").Value = 10.94309780184223
$ws.Range("B24").Value = 21.06630964485228
$ws.Range("C24").Value = 11.5383631584292
$ws.Range("D24").Value = 6.877313818657475
$ws.Range("F24").Value = 39.45971360661974
$ws.Range("G24").Value = 3.666693582022686
$ws.Range("I24").Value = 30.43807698682037
$ws.Range("L24").Value = 10.87811868649801
$ws.Range("B25").Value = 20.15788523278878
$ws.Range("C25").Value = 10.47906143751753
$ws.Range("D25").Value = 6.905569780953721
$ws.Range("F25").Value = 38.58023818262608
$ws.Range("G25").Value = 3.674203650794795
$ws.Range("I25").Value = 30.16402329294877
$ws.Range("L25").Value = 10.81668839572379
